# Updates the cryptos price/volume table with the latest scrape.
# Column D ("Price") values are numeric-looking text (e.g. "0.4830",
# "27.847.45") that must stay stored as text, exactly as scraped, so we
# force NumberFormat "@" (Text) before writing and then reset the style
# back to "Normal" so we don't leave a stray per-cell format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.847.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4830"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3816"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9399"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07796"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.503"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.611"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008865"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.858.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.119"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.107.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.948"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.053"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.983"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08896"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.336"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.231"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.736"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.129"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02044"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5618"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.057"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.556"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1529"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4880"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06123"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
